$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    2.4218607543199702,
    2.4168409172640701,
    2.4121668977944499,
    2.4100298601538701,
    2.4116420813322001,
    2.4172568023318801,
    2.4373918210577599,
    2.4697596502624002,
    2.5141890136612299,
    2.5689125726683302,
    2.6303681608079899,
    2.6850154328682798,
    2.7278158539492598,
    2.7606990865499301,
    2.7860610063703799,
    2.8064014499033898,
    2.8242653465606899,
    2.85067910568225,
    2.87849166622525,
    2.9020316800600598,
    2.9188738991237502,
    2.9088603104899802,
    2.8873433960221799,
    2.8616848819179701,
    2.83818153745066,
    2.8225164024903502,
    2.8158776190786199,
    2.8054602400362199,
    2.7894702641117699,
    2.76796424320393,
    2.7426066411610499,
    2.7166427988810198,
    2.7000229584351501,
    2.6914320833608198,
    2.6876725281839402,
    2.6842333833533099,
    2.6842837430066999,
    2.6851080643265202,
    2.6885119022639001,
    2.69718625631387,
    2.7159852335500001,
    2.7292375648642602,
    2.7338451603668101,
    2.73210460856779,
    2.7265288244901802,
    2.7158629240487802,
    2.7131151090771901,
    2.7119162316099699,
    2.7120455972025299,
    2.71460634618896,
    2.7224446088099001,
    2.7302447334506299,
    2.7415327372310601,
    2.7552677554793399,
    2.7689982341447301,
    2.7785605320836102,
    2.7888846398579799,
    2.7997698458821501,
    2.8127048231521301,
    2.8284917281343498,
    2.8410258815829699,
    2.8584153349424999,
    2.8792468946735199,
    2.9036101646441099,
    2.93041212925089,
    2.9470324878774599,
    2.9575658532107898,
    2.96694648628847,
    2.9756864927511502,
    2.9838340320744599,
    3.00148395640409,
    3.0232375660370598,
    3.04498283190396,
    3.0622650416916501,
    3.0712693832640001,
    3.0788860611276401,
    3.0827868163175798,
    3.0886368479509101,
    3.1009995338636802,
    3.1232052557622199,
    3.1465207006423102,
    3.1708111710055298,
    3.1907830300968199,
    3.20547305719812,
    3.2162633618004799,
    3.2182495648441001,
    3.21009886538255,
    3.1919670033640402,
    3.16129521755416,
    3.1141591729089901,
    3.05115811498835,
    2.9671967791762301,
    2.87926616388215,
    2.7970376200081999,
    2.7251935022059799,
    2.6773484811531501,
    2.6453865401352599,
    2.6171288754555899,
    2.5877200226106098,
    2.5587272476943999
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}
